$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-6 from
# serial 45175 (2023-09-06) to serial 45183 (2023-09-14)
$newDate = 45183

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
$ws.Range("C6").Value = $newDate
